# Update projections + add jupyter
$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("UpdatedResults")
$ws2 = $wb.Worksheets.Item("Archive")
$ws3 = $wb.Worksheets.Item("WL Record")

# ---------------------------------------------------------------------------
# 1) UpdatedResults sheet: the old prediction rows (2-4) have now been
#    played out and move down to the Archive sheet (below), and this sheet
#    is refilled with the *new* set of upcoming-game predictions.
# ---------------------------------------------------------------------------

$sheet1Rows = @(
    @{ A=3;    B="Cleveland +3";              C="Washington Wizards";     D=102; E="Cleveland Cavaliers";    F=107; G=209; H=211.5; I="UNDER" },
    @{ A=9;    B="Orlando +9";                C="Brooklyn Nets";          D=101; E="Orlando Magic";          F=104; G=205; H=209.5; I="UNDER" },
    @{ A=-2;   B="Boston - 2";                C="Toronto Raptors";        D=102; E="Boston Celtics";         F=105; G=207; H=209.5; I="UNDER" },
    @{ A=-3.5; B="Houston -3.5";              C="Detroit Pistons";        D=95;  E="Houston Rockets";        F=106; G=201; H=210.5; I="UNDER" },
    @{ A=-3;   B="New York -3";               C="Milwaukee Bucks";        D=106; E="New York Knicks";        F=114; G=220; H=215.5; I="OVER"  },
    @{ A=-3.5; B="Chicago -3";                C="Dallas Mavericks";       D=100; E="Chicago Bulls";          F=107; G=207; H=213;   I="UNDER" },
    @{ A=-4.5; B="Charlotte +4.5";            C="Charlotte Hornets";      D=117; E="Memphis Grizzlies";      F=115; G=232; H=228.5; I="OVER"  },
    @{ A=-4.5; B="New Orleans -4.5";          C="Oklahoma City Thunder";  D=97;  E="New Orleans Pelicans";   F=108; G=205; H=211;   I="UNDER" },
    @{ A=0;    B="Sacramento - Pick";         C="Sacramento Kings";       D=109; E="San Antonio Spurs";      F=104; G=213; H=220;   I="UNDER" },
    @{ A=3;    B="Indiana - 3";               C="Indiana Pacers";         D=107; E="Denver Nuggets";         F=96;  G=203; H=213;   I="UNDER" },
    @{ A=-6;   B="Portland +6";               C="Portland Trail Blazers"; D=110; E="Phoenix Suns";           F=111; G=221; H=220.5; I="OVER"  },
    @{ A=-7;   B="Minnesota +7";              C="Minnesota Timberwolves"; D=101; E="Golden State Warriors";  F=106; G=207; H=222.5; I="UNDER" },
    @{ A=4.5;  B="Los Angeles Lakers +4.5";   C="Miami Heat";             D=104; E="Los Angeles Lakers";     F=110; G=214; H=214.5; I="UNDER" }
)

$r = 2
foreach ($row in $sheet1Rows) {
    $ws1.Cells.Item($r, 1).Value = $row.A
    $ws1.Cells.Item($r, 2).Value = $row.B
    $ws1.Cells.Item($r, 3).Value = $row.C
    $ws1.Cells.Item($r, 4).Value = $row.D
    $ws1.Cells.Item($r, 5).Value = $row.E
    $ws1.Cells.Item($r, 6).Value = $row.F
    $ws1.Cells.Item($r, 7).Value = $row.G
    $ws1.Cells.Item($r, 8).Value = $row.H
    $ws1.Cells.Item($r, 9).Value = $row.I
    $r++
}

# ---------------------------------------------------------------------------
# 2) Archive sheet: append the 3 now-completed games (previously rows 2-4
#    of UpdatedResults) as final rows 49-51, copying the row-18 number
#    formatting/style (A:Q) so the new rows match the existing block.
# ---------------------------------------------------------------------------

$ws2.Range("A18:Q18").Copy() | Out-Null
$ws2.Range("A49:Q51").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$archiveRows = @(
    @{ Row=49; A=44509; B="Milwaukee Bucks";        C=118; D="Philadelphia 76ers";   E=109; G=6.5; H="LOSS"; J=105; L=112; N=220.5 },
    @{ Row=50; A=44509; B="Atlanta Hawks";           C=98;  D="Utah Jazz";            E=110; G=-8;  H="LOSS"; J=112; L=110; N=221.5 },
    @{ Row=51; A=44509; B="Portland Trail Blazers";  C=109; D="Los Angeles Clippers"; E=117; G=-3;  H="LOSS"; J=109; L=102; N=220   }
)

foreach ($row in $archiveRows) {
    $rr = $row.Row
    $ws2.Cells.Item($rr, 1).Value = $row.A
    $ws2.Cells.Item($rr, 2).Value = $row.B
    $ws2.Cells.Item($rr, 3).Value = $row.C
    $ws2.Cells.Item($rr, 4).Value = $row.D
    $ws2.Cells.Item($rr, 5).Value = $row.E
    $ws2.Cells.Item($rr, 6).Formula = "=SUM(J$rr - L$rr)"
    $ws2.Cells.Item($rr, 7).Value = $row.G
    $ws2.Cells.Item($rr, 8).Value = $row.H
    $ws2.Cells.Item($rr, 9).Value = $row.B
    $ws2.Cells.Item($rr, 10).Value = $row.J
    $ws2.Cells.Item($rr, 11).Value = $row.D
    $ws2.Cells.Item($rr, 12).Value = $row.L
    $ws2.Cells.Item($rr, 13).Formula = "=SUM(J$rr+L$rr)"
    $ws2.Cells.Item($rr, 14).Value = $row.N
    $ws2.Cells.Item($rr, 15).Formula = "=SUM(C$rr+E$rr)"
    $ws2.Cells.Item($rr, 16).Formula = "=IF(M$rr<N$rr,""UNDER"",""OVER"")"
    $ws2.Cells.Item($rr, 17).Formula = "=IF(O$rr<N$rr,""UNDER"",""OVER"")"
    $ws2.Cells.Item($rr, 18).Formula = "=IF(P$rr=Q$rr,""WIN"",""LOSS"")"
}

# ---------------------------------------------------------------------------
# 3) View / selection state
# ---------------------------------------------------------------------------

$ws2.Columns.Item(4).ColumnWidth = 20.5

$ws1.Activate() | Out-Null
$ws1.Range("D14").Select() | Out-Null

$ws2.Activate() | Out-Null
$ws2.Range("I42").Select() | Out-Null

$ws3.Activate() | Out-Null
$ws3.Range("F2").Select() | Out-Null

Write-Host "Applied NBA projection update"
